$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02942666666666667
$ws.Range("H2").Value = 0.08828
$ws.Range("I2").Value = 0.02893416853923747
$ws.Range("J2").Value = 0.02893416853923746
$ws.Range("M2").Value = 0.254776
$ws.Range("N2").Value = 0.764328
$ws.Range("O2").Value = 0.2586256426914624
$ws.Range("P2").Value = 0.2586256426914624
$ws.Range("Q2").Value = 0.007497208426666666
$ws.Range("R2").Value = 0.06747487584
$ws.Range("S2").Value = 0.007483117934203382
$ws.Range("T2").Value = 0.00748311793420338
$ws.Range("G3").Value = 0.02942666666666667
$ws.Range("H3").Value = 0.08828
$ws.Range("I3").Value = 0.02893416853923747
$ws.Range("J3").Value = 0.02893416853923746
$ws.Range("O3").Value = 0.4481730559376317
$ws.Range("P3").Value = 0.4481730559376316
$ws.Range("Q3").Value = 0.01299193218666667
$ws.Range("R3").Value = 0.11692738968
$ws.Range("S3").Value = 0.01296751473524454
$ws.Range("T3").Value = 0.01296751473524453
$ws.Range("G4").Value = 0.02942666666666667
$ws.Range("H4").Value = 0.08828
$ws.Range("I4").Value = 0.02893416853923747
$ws.Range("J4").Value = 0.02893416853923746
$ws.Range("M4").Value = 0.288837
$ws.Range("N4").Value = 0.866511
$ws.Range("O4").Value = 0.2932013013709059
$ws.Range("P4").Value = 0.2932013013709059
$ws.Range("Q4").Value = 0.00849951012
$ws.Range("R4").Value = 0.07649559107999999
$ws.Range("S4").Value = 0.00848353586978955
$ws.Range("T4").Value = 0.008483535869789546
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04158733333333333
$ws.Range("H5").Value = 0.124762
$ws.Range("I5").Value = 0.04089130873688655
$ws.Range("J5").Value = 0.04089130873688654
$ws.Range("M5").Value = 0.254776
$ws.Range("N5").Value = 0.764328
$ws.Range("O5").Value = 0.2586256426914624
$ws.Range("P5").Value = 0.2586256426914624
$ws.Range("Q5").Value = 0.01059545443733333
$ws.Range("R5").Value = 0.09535908993599999
$ws.Range("S5").Value = 0.0105755410025723
$ws.Range("T5").Value = 0.01057554100257229
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04158733333333333
$ws.Range("H6").Value = 0.124762
$ws.Range("I6").Value = 0.04089130873688655
$ws.Range("J6").Value = 0.04089130873688654
$ws.Range("O6").Value = 0.4481730559376317
$ws.Range("P6").Value = 0.4481730559376316
$ws.Range("Q6").Value = 0.01836089084133333
$ws.Range("R6").Value = 0.165248017572
$ws.Range("S6").Value = 0.01832638279789962
$ws.Range("T6").Value = 0.01832638279789962
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04158733333333333
$ws.Range("H7").Value = 0.124762
$ws.Range("I7").Value = 0.04089130873688655
$ws.Range("J7").Value = 0.04089130873688654
$ws.Range("M7").Value = 0.288837
$ws.Range("N7").Value = 0.866511
$ws.Range("O7").Value = 0.2932013013709059
$ws.Range("P7").Value = 0.2932013013709059
$ws.Range("Q7").Value = 0.012011960598
$ws.Range("R7").Value = 0.108107645382
$ws.Range("S7").Value = 0.01198938493641463
$ws.Range("T7").Value = 0.01198938493641463
$ws.Range("G8").Value = 0.9379063333333333
$ws.Range("H8").Value = 2.813719
$ws.Range("I8").Value = 0.9222091047582091
$ws.Range("J8").Value = 0.9222091047582088
$ws.Range("M8").Value = 0.254776
$ws.Range("N8").Value = 0.764328
$ws.Range("O8").Value = 0.2586256426914624
$ws.Range("P8").Value = 0.2586256426914624
$ws.Range("Q8").Value = 0.2389560239813333
$ws.Range("R8").Value = 2.150604215832
$ws.Range("S8").Value = 0.23850692241401
$ws.Range("T8").Value = 0.23850692241401
$ws.Range("G9").Value = 0.9379063333333333
$ws.Range("H9").Value = 2.813719
$ws.Range("I9").Value = 0.9222091047582091
$ws.Range("J9").Value = 0.9222091047582088
$ws.Range("O9").Value = 0.4481730559376317
$ws.Range("P9").Value = 0.4481730559376316
$ws.Range("Q9").Value = 0.4140875219793333
$ws.Range("R9").Value = 3.726787697814
$ws.Range("S9").Value = 0.4133092726929941
$ws.Range("T9").Value = 0.4133092726929939
$ws.Range("G10").Value = 0.9379063333333333
$ws.Range("H10").Value = 2.813719
$ws.Range("I10").Value = 0.9222091047582091
$ws.Range("J10").Value = 0.9222091047582088
$ws.Range("M10").Value = 0.288837
$ws.Range("N10").Value = 0.866511
$ws.Range("O10").Value = 0.2932013013709059
$ws.Range("P10").Value = 0.2932013013709059
$ws.Range("Q10").Value = 0.270902051601
$ws.Range("R10").Value = 2.438118464409
$ws.Range("S10").Value = 0.270392909651205
$ws.Range("T10").Value = 0.2703929096512049
$ws.Range("G11").Value = 0.008101000000000001
$ws.Range("H11").Value = 0.024303
$ws.Range("I11").Value = 0.007965417965667061
$ws.Range("J11").Value = 0.007965417965667059
$ws.Range("M11").Value = 0.254776
$ws.Range("N11").Value = 0.764328
$ws.Range("O11").Value = 0.2586256426914624
$ws.Range("P11").Value = 0.2586256426914624
$ws.Range("Q11").Value = 0.002063940376
$ws.Range("R11").Value = 0.018575463384
$ws.Range("S11").Value = 0.002060061340676765
$ws.Range("T11").Value = 0.002060061340676764
$ws.Range("G12").Value = 0.008101000000000001
$ws.Range("H12").Value = 0.024303
$ws.Range("I12").Value = 0.007965417965667061
$ws.Range("J12").Value = 0.007965417965667059
$ws.Range("O12").Value = 0.4481730559376317
$ws.Range("P12").Value = 0.4481730559376316
$ws.Range("Q12").Value = 0.003576607702
$ws.Range("R12").Value = 0.032189469318
$ws.Range("S12").Value = 0.00356988571149352
$ws.Range("T12").Value = 0.003569885711493519
$ws.Range("G13").Value = 0.008101000000000001
$ws.Range("H13").Value = 0.024303
$ws.Range("I13").Value = 0.007965417965667061
$ws.Range("J13").Value = 0.007965417965667059
$ws.Range("M13").Value = 0.288837
$ws.Range("N13").Value = 0.866511
$ws.Range("O13").Value = 0.2932013013709059
$ws.Range("P13").Value = 0.2932013013709059
$ws.Range("Q13").Value = 0.002339868537
$ws.Range("R13").Value = 0.021058816833
$ws.Range("S13").Value = 0.002335470913496777
$ws.Range("T13").Value = 0.002335470913496776
